$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("496").Delete()
